$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "Compression ratios" block (title + Tabela5 table, columns E:H)
# moves up by one row: old title row 3 / table E4:H7  ->
# new title row 2 / table E3:H6. We rebuild it explicitly (cell by
# cell, with formats applied directly) rather than relying on
# row-shifting Delete/Insert, because a partial-column Delete/Insert
# in this host shifts the WHOLE row (including columns B:C), which
# must stay exactly where it is.
# ------------------------------------------------------------------

$xlCenter = -4108
$xlUnderlineStyleSingle = 2
$xlUnderlineStyleNone = -4142
$numFmtAccounting = '_-* #,##0.00_-;-* #,##0.00_-;_-* "-"??_-;_-@_-'

# 1) Remove the old merged title and clear the whole old block E3:H7.
$ws.Range("E3:H3").UnMerge()
$ws.Range("E3:H7").Clear()

# 2) Title row (row 2), merged E2:H2 - same look as B2/B9/B16 titles
#    (horizontal-center only, default vertical alignment).
$ws.Range("E2:H2").Merge()
$ws.Range("E2").Value = "Compression ratios"
$ws.Range("E2:H2").HorizontalAlignment = $xlCenter

# 3) Header row (row 3): quality / barn / logo / peppers
$ws.Range("E3").Value = "quality"
$ws.Range("F3").Value = "barn"
$ws.Range("G3").Value = "logo"
$ws.Range("H3").Value = "peppers"
$ws.Range("E3:H3").HorizontalAlignment = $xlCenter
$ws.Range("E3:H3").VerticalAlignment = $xlCenter
$ws.Range("H3").Font.Underline = $xlUnderlineStyleSingle

# 4) Data rows 4-6: high/medium/low plus ratio formulas.
$ws.Range("E4").Value = "high"
$ws.Range("F4").Formula = "=ROUND(`$C`$4/C5,0)&`":`"&1"
$ws.Range("G4").Formula = "=ROUND(`$C`$11/C12,0)&`":`"&1"
$ws.Range("H4").Formula = "=ROUND(`$C`$18/C19,0)&`":`"&1"

$ws.Range("E5").Value = "medium"
$ws.Range("F5").Formula = "=ROUND(`$C`$4/C6,0)&`":`"&1"
$ws.Range("G5").Formula = "=ROUND(`$C`$11/C13,0)&`":`"&1"
$ws.Range("H5").Formula = "=ROUND(`$C`$18/C20,0)&`":`"&1"

$ws.Range("E6").Value = "low"
$ws.Range("F6").Formula = "=ROUND(`$C`$4/C7,0)&`":`"&1"
$ws.Range("G6").Formula = "=ROUND(`$C`$11/C14,0)&`":`"&1"
$ws.Range("H6").Formula = "=ROUND(`$C`$18/C21,0)&`":`"&1"

$ws.Range("E4:E6").HorizontalAlignment = $xlCenter
$ws.Range("E4:E6").VerticalAlignment = $xlCenter

$ws.Range("F4:H6").HorizontalAlignment = $xlCenter
$ws.Range("F4:H6").VerticalAlignment = $xlCenter
$ws.Range("F4:H6").NumberFormat = $numFmtAccounting

# 5) Resize the Tabela5 list object onto its new location and flip on
#    "first column" emphasis (matches the other three tables).
$lo = $ws.ListObjects.Item("Tabela5")
$lo.Resize($ws.Range("E3:H6"))
$lo.ShowTableStyleFirstColumn = $true

# 6) Drop the stray formatted-but-empty M8 cell.
$ws.Range("M8").Clear()

# 7) Cursor/selection cosmetics.
$ws.Range("L11").Select()
